$d = $word.ActiveDocument

# Target fill colors: change old gray (B3B3B3) shading to new gray (CCCCCC)
# on table cells, wherever it occurs in the document's tables.
$oldFill = 0xB3B3B3
$newFill = 0xCCCCCC

foreach ($tbl in $d.Tables) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $cell = $tbl.Cell($r, $c)
            if ($cell.Shading.BackgroundPatternColor -eq $oldFill) {
                $cell.Shading.BackgroundPatternColor = $newFill
            }
        }
    }
}
